# Add new simulation results (columns V:AO) to Sheet1, mirroring the existing
# B:U layout: row 2 holds the angle samples (0..2*pi over 20 steps) and row 3
# holds the corresponding 0/1 "loaded" flags for the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (angles): identical 20-value sequence already present in B2:U2.
$row2Values = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3 (flags): new simulation outcomes for columns V:AO.
$row3Values = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 0, 1, 1, 1)

# Columns V (22) through AO (41).
$startCol = 22
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Scroll the view toward the newly added columns and move the selection,
# matching the author's on-screen state after adding the data.
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("U8").Select() | Out-Null
